$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "https://dev.recordati-plus.de/de_DE/overview-page/home-psychiatrie",
    "https://dev.recordati-plus.de/de_DE/overview-page/home-urologie",
    "https://dev.recordati-plus.de/de_DE/overview-page/alle-inhalte-gastroenterologie",
    "https://dev.recordati-plus.de/de_DE/overview-page/alle-inhalte-psychiatrie",
    "https://dev.recordati-plus.de/de_DE/overview-page/alle-inhalte-paediatrie",
    "https://dev.recordati-plus.de/de_DE/overview-page/alle-inhalte-kardiologie",
    "https://dev.recordati-plus.de/de_DE/overview-page/alle-inhalte-orthopaedie",
    "https://dev.recordati-plus.de/de_DE/overview-page/alle-inhalte-urologie",
    "https://dev.recordati-plus.de/de_DE/overview-page/alle-inhalte-erkaeltung",
    "https://dev.recordati-plus.de/de_DE/overview-page/alle-inhalte-sonstige-praeparate"
)

$row = 3
foreach ($url in $urls) {
    $ws.Cells.Item($row, 1).Value = $url
    $row++
}
